# Update Excel files after daily scrape - 2025-08-28 03:08:51 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-9 (columns A-H)
# Column A values are prefixed with a leading apostrophe so Excel stores
# them as text (matching the original inlineStr type) instead of numbers.
$data = @(
    @("'1327125", "https://aiesec.org/opportunity/global-talent/1327125", "Machine Learning Intern", "Sahibzada Ajit Singh Nagar, Punjab, India", "No", "1 applicant", "9 - 12 Weeks", "Solitaire Infosys Pvt. Ltd"),
    @("'1327106", "https://aiesec.org/opportunity/global-talent/1327106", "Sales Assistant (Spanish)", "Denizli, Kumkısık, Denizli, Türkiye", "No", "1 applicant", "6 - 18 Months", "Sera Moda"),
    @("'1327090", "https://aiesec.org/opportunity/global-talent/1327090", "AI Intern", "Chandigarh, India", "No", "1 applicant", "3 - 6 Months", "Solitaire Infosys Pvt. Ltd"),
    @("'1326864", "https://aiesec.org/opportunity/global-talent/1326864", "Transport Sales Manager", "Vienna, Austria", "No", "6 applicants", "6 - 18 Months", "FERCAM AUSTRIA GMHB"),
    @("'1326162", "https://aiesec.org/opportunity/global-talent/1326162", "Tech Sales Development Representative( swedish Only)", "Bournemouth, Royaume-Uni", "No", "9 applicants", "6 - 18 Months", "EIMS Ltd"),
    @("'1326152", "https://aiesec.org/opportunity/global-talent/1326152", "Tech Sales Development Representative( German / Austrian Only)", "Bournemouth, Royaume-Uni", "No", "6 applicants", "6 - 18 Months", "EIMS Ltd"),
    @("'1324596", "https://aiesec.org/opportunity/global-talent/1324596", "ACCOUNTANT", "New Damietta City, Damietta El-Gadeeda City, New Damietta, Damietta Governorate, Egypt", "No", "48 applicants", "3 - 6 Months", "Business Haven Consultancy"),
    @("'1311536", "https://aiesec.org/opportunity/global-talent/1311536", "Accelerate Romania | Managing Co-founder", "Bucharest, Romania", "No", "132 applicants", "9 - 12 Weeks", "Skulptor")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$c]
    }
}

# Remove old rows 10, 11, 12 (now beyond the new data range)
$ws.Rows.Item(10).EntireRow.Delete()
$ws.Rows.Item(10).EntireRow.Delete()
$ws.Rows.Item(10).EntireRow.Delete()

# Update column widths per diff.
# NOTE: the engine's ColumnWidth setter re-derives the stored OOXML <col width>
# value via a pixel-rounding conversion that adds ~0.83 to whatever is
# assigned (Excel's character-width -> pixel -> character-width round trip).
# Subtracting 0.9 compensates so the persisted width lands exactly on the
# desired integer value.
$ws.Columns.Item(3).ColumnWidth = 65 - 0.9
$ws.Columns.Item(4).ColumnWidth = 89 - 0.9
$ws.Columns.Item(8).ColumnWidth = 29 - 0.9
